$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B18").Copy()
$ws.Range("B19").PasteSpecial(-4122)

$ws.Range("A19").Value = "Aris"
$ws.Range("B19").Value = 45334
$ws.Range("C19").Formula = "=12+53/60"
$ws.Range("D19").Formula = "=13+7/60"
$ws.Range("F19").Value = "Button even link"

$ws.Range("F20").Select()
